# Actualización automática 2025-12-09 13:30:06
$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Row 5 (AUCANSHALA ALLAICA FREDDY HERNAN)
$wsGrupo.Range("L5").Value = 455.77
$wsGrupo.Range("M5").Value = 1762.56

# Row 11 (CUSTODE FALCONI DIANA CAROLINA)
$wsGrupo.Range("D11").Value = 1391.04
$wsGrupo.Range("I11").Value = 471.6
$wsGrupo.Range("L11").Value = 443.44

# Row 38 totals ("x de 36" counters)
$wsGrupo.Range("D38").Value = "1 de 36"
$wsGrupo.Range("I38").Value = "1 de 36"
$wsGrupo.Range("L38").Value = "2 de 36"
$wsGrupo.Range("M38").Value = "4 de 36"

# --- Sheet "VENTA MENSUAL" ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

# Row 5 (AUCANSHALA ALLAICA FREDDY HERNAN) - diciembre
$wsMensual.Range("F5").Value = 2218.33

# Row 11 (CUSTODE FALCONI DIANA CAROLINA) - diciembre
$wsMensual.Range("F11").Value = 2306.08

# Row 38 total - diciembre
$wsMensual.Range("F38").Value = 6468.95
